# Apply numeric updates to match the target OOXML diff.
# All cells in this workbook are static values (no formulas),
# so every changed cell is written explicitly per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0.05227602163694004
$ws.Range("E2").Value = 0.1446971150797848
$ws.Range("F2").Value = 0.2159902207704085
$ws.Range("G2").Value = 0.2172915037649183
$ws.Range("H2").Value = 0.001261706832069276
$ws.Range("I2").Value = 0.04225580596258013
$ws.Range("J2").Value = 0.02002500917777559
$ws.Range("N2").Value = 0.02383260028749569
$ws.Range("O2").Value = 0.06408475540338457
$ws.Range("P2").Value = 0.02553444209348987
$ws.Range("R2").Value = 0.03572876697831024
$ws.Range("S2").Value = 0.07227525240818347
$ws.Range("AA2").Value = 0.004338494239685599
$ws.Range("AB2").Value = 0.04104287734301845
$ws.Range("AD2").Value = 0.03805326964814988
$ws.Range("AG2").Value = 0.001312158373805273
$ws.Range("F3").Value = 0.2963964809008053
$ws.Range("G3").Value = 0.002965733538903559
$ws.Range("H3").Value = 0.2650021891096919
$ws.Range("I3").Value = 0.0096831372908373
$ws.Range("J3").Value = 0.04613074204241689
$ws.Range("L3").Value = 0.06836262997400276
$ws.Range("M3").Value = 0.02661474433962446
$ws.Range("P3").Value = 0.01590951728843793
$ws.Range("Q3").Value = 0.09341925912490752
$ws.Range("S3").Value = 0.05009577398566298
$ws.Range("T3").Value = 0.03712953881284548
$ws.Range("AC3").Value = 0.0337780559748183
$ws.Range("AD3").Value = 0.002147524538128404
$ws.Range("AE3").Value = 0.05236467307891722
$ws.Range("E4").Value = 0.2690853829751593
$ws.Range("G4").Value = 0.1730596773045425
$ws.Range("H4").Value = 0.03410949806375992
$ws.Range("I4").Value = 0.009864619963309946
$ws.Range("K4").Value = 0.08866487474151809
$ws.Range("M4").Value = 0.1307204712476435
$ws.Range("P4").Value = 0.03614280509614683
$ws.Range("Q4").Value = 0.07851399728811799
$ws.Range("R4").Value = 0.003088005375383761
$ws.Range("S4").Value = 0.1180369118025626
$ws.Range("Z4").Value = 0.004594321142818709
$ws.Range("AC4").Value = 0.0312059318817504
$ws.Range("AE4").Value = 0.02086900131924758
$ws.Range("AH4").Value = 0.002044501798038871
$ws.Range("E5").Value = 0.1135939084770893
$ws.Range("F5").Value = 0.03799156432914786
$ws.Range("G5").Value = 0.2321456331076833
$ws.Range("H5").Value = 0.1507498863031252
$ws.Range("I5").Value = 0.01125732688188755
$ws.Range("J5").Value = 0.04587723851279347
$ws.Range("K5").Value = 0.05332320784153089
$ws.Range("M5").Value = 0.02243360034225561
$ws.Range("O5").Value = 0.0004771827899888273
$ws.Range("P5").Value = 0.1028473450473246
$ws.Range("Q5").Value = 0.05635548472940122
$ws.Range("S5").Value = 0.07157356726276075
$ws.Range("T5").Value = 0.009750854928658364
$ws.Range("Z5").Value = 0.001898044066896292
$ws.Range("AC5").Value = 0.05628988292730669
$ws.Range("AE5").Value = 0.02663998569534371
$ws.Range("AJ5").Value = 0.006795286756806496
$ws.Range("D6").Value = 0.02907864705958124
$ws.Range("E6").Value = 0.2592514713256358
$ws.Range("G6").Value = 0.1616800671913459
$ws.Range("H6").Value = 0.04974572538659951
$ws.Range("K6").Value = 0.04162897780336663
$ws.Range("L6").Value = 0.02164350736664281
$ws.Range("M6").Value = 0.07289734555966348
$ws.Range("N6").Value = 0.03476709252170095
$ws.Range("O6").Value = 0.01165902253182157
$ws.Range("P6").Value = 0.008944451578672482
$ws.Range("Q6").Value = 0.1304176016924338
$ws.Range("S6").Value = 0.143425360499828
$ws.Range("T6").Value = 0.008036127199859631
$ws.Range("Y6").Value = 0.00223209534457359
$ws.Range("AA6").Value = 0.01316612998742762
$ws.Range("AC6").Value = 0.01099808812383458
$ws.Range("AE6").Value = 0.0004282888270125547

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0.05227602163694004
$ws.Range("E2").Value = 0.1969731367167249
$ws.Range("F2").Value = 0.4129633574871333
$ws.Range("G2").Value = 0.6302548612520517
$ws.Range("H2").Value = 0.6315165680841209
$ws.Range("I2").Value = 0.6737723740467011
$ws.Range("J2").Value = 0.6937973832244767
$ws.Range("K2").Value = 0.6937973832244767
$ws.Range("L2").Value = 0.6937973832244767
$ws.Range("M2").Value = 0.6937973832244767
$ws.Range("N2").Value = 0.7176299835119724
$ws.Range("O2").Value = 0.7817147389153569
$ws.Range("P2").Value = 0.8072491810088468
$ws.Range("Q2").Value = 0.8072491810088468
$ws.Range("R2").Value = 0.8429779479871571
$ws.Range("S2").Value = 0.9152532003953405
$ws.Range("T2").Value = 0.9152532003953405
$ws.Range("U2").Value = 0.9152532003953405
$ws.Range("V2").Value = 0.9152532003953405
$ws.Range("W2").Value = 0.9152532003953405
$ws.Range("X2").Value = 0.9152532003953405
$ws.Range("Y2").Value = 0.9152532003953405
$ws.Range("Z2").Value = 0.9152532003953405
$ws.Range("AA2").Value = 0.9195916946350261
$ws.Range("AB2").Value = 0.9606345719780446
$ws.Range("AC2").Value = 0.9606345719780446
$ws.Range("AD2").Value = 0.9986878416261945
$ws.Range("AE2").Value = 0.9986878416261945
$ws.Range("AF2").Value = 0.9986878416261945
$ws.Range("F3").Value = 0.2963964809008053
$ws.Range("G3").Value = 0.2993622144397089
$ws.Range("H3").Value = 0.5643644035494009
$ws.Range("I3").Value = 0.5740475408402381
$ws.Range("J3").Value = 0.620178282882655
$ws.Range("K3").Value = 0.620178282882655
$ws.Range("L3").Value = 0.6885409128566577
$ws.Range("M3").Value = 0.7151556571962822
$ws.Range("N3").Value = 0.7151556571962822
$ws.Range("O3").Value = 0.7151556571962822
$ws.Range("P3").Value = 0.7310651744847201
$ws.Range("Q3").Value = 0.8244844336096276
$ws.Range("R3").Value = 0.8244844336096276
$ws.Range("S3").Value = 0.8745802075952906
$ws.Range("T3").Value = 0.9117097464081362
$ws.Range("U3").Value = 0.9117097464081362
$ws.Range("V3").Value = 0.9117097464081362
$ws.Range("W3").Value = 0.9117097464081362
$ws.Range("X3").Value = 0.9117097464081362
$ws.Range("Y3").Value = 0.9117097464081362
$ws.Range("Z3").Value = 0.9117097464081362
$ws.Range("AA3").Value = 0.9117097464081362
$ws.Range("AB3").Value = 0.9117097464081362
$ws.Range("AC3").Value = 0.9454878023829545
$ws.Range("AD3").Value = 0.9476353269210829
$ws.Range("AE3").Value = 1
$ws.Range("AF3").Value = 1
$ws.Range("AG3").Value = 1
$ws.Range("AH3").Value = 1
$ws.Range("AI3").Value = 1
$ws.Range("AJ3").Value = 1
$ws.Range("AK3").Value = 1
$ws.Range("E4").Value = 0.2690853829751593
$ws.Range("F4").Value = 0.2690853829751593
$ws.Range("G4").Value = 0.4421450602797018
$ws.Range("H4").Value = 0.4762545583434617
$ws.Range("I4").Value = 0.4861191783067716
$ws.Range("J4").Value = 0.4861191783067716
$ws.Range("K4").Value = 0.5747840530482897
$ws.Range("L4").Value = 0.5747840530482897
$ws.Range("M4").Value = 0.7055045242959332
$ws.Range("N4").Value = 0.7055045242959332
$ws.Range("O4").Value = 0.7055045242959332
$ws.Range("P4").Value = 0.7416473293920801
$ws.Range("Q4").Value = 0.8201613266801981
$ws.Range("R4").Value = 0.8232493320555819
$ws.Range("S4").Value = 0.9412862438581445
$ws.Range("T4").Value = 0.9412862438581445
$ws.Range("U4").Value = 0.9412862438581445
$ws.Range("V4").Value = 0.9412862438581445
$ws.Range("W4").Value = 0.9412862438581445
$ws.Range("X4").Value = 0.9412862438581445
$ws.Range("Y4").Value = 0.9412862438581445
$ws.Range("Z4").Value = 0.9458805650009632
$ws.Range("AA4").Value = 0.9458805650009632
$ws.Range("AB4").Value = 0.9458805650009632
$ws.Range("AC4").Value = 0.9770864968827135
$ws.Range("AD4").Value = 0.9770864968827135
$ws.Range("AE4").Value = 0.9979554982019612
$ws.Range("AF4").Value = 0.9979554982019612
$ws.Range("AG4").Value = 0.9979554982019612
$ws.Range("E5").Value = 0.1135939084770893
$ws.Range("F5").Value = 0.1515854728062371
$ws.Range("G5").Value = 0.3837311059139205
$ws.Range("H5").Value = 0.5344809922170457
$ws.Range("I5").Value = 0.5457383190989332
$ws.Range("J5").Value = 0.5916155576117267
$ws.Range("K5").Value = 0.6449387654532576
$ws.Range("L5").Value = 0.6449387654532576
$ws.Range("M5").Value = 0.6673723657955132
$ws.Range("N5").Value = 0.6673723657955132
$ws.Range("O5").Value = 0.6678495485855021
$ws.Range("P5").Value = 0.7706968936328267
$ws.Range("Q5").Value = 0.8270523783622279
$ws.Range("R5").Value = 0.8270523783622279
$ws.Range("S5").Value = 0.8986259456249887
$ws.Range("T5").Value = 0.908376800553647
$ws.Range("U5").Value = 0.908376800553647
$ws.Range("V5").Value = 0.908376800553647
$ws.Range("W5").Value = 0.908376800553647
$ws.Range("X5").Value = 0.908376800553647
$ws.Range("Y5").Value = 0.908376800553647
$ws.Range("Z5").Value = 0.9102748446205433
$ws.Range("AA5").Value = 0.9102748446205433
$ws.Range("AB5").Value = 0.9102748446205433
$ws.Range("AC5").Value = 0.9665647275478499
$ws.Range("AD5").Value = 0.9665647275478499
$ws.Range("AE5").Value = 0.9932047132431937
$ws.Range("AF5").Value = 0.9932047132431937
$ws.Range("AG5").Value = 0.9932047132431937
$ws.Range("AH5").Value = 0.9932047132431937
$ws.Range("AI5").Value = 0.9932047132431937
$ws.Range("D6").Value = 0.02907864705958124
$ws.Range("E6").Value = 0.288330118385217
$ws.Range("F6").Value = 0.288330118385217
$ws.Range("G6").Value = 0.4500101855765629
$ws.Range("H6").Value = 0.4997559109631624
$ws.Range("I6").Value = 0.4997559109631624
$ws.Range("J6").Value = 0.4997559109631624
$ws.Range("K6").Value = 0.541384888766529
$ws.Range("L6").Value = 0.5630283961331718
$ws.Range("M6").Value = 0.6359257416928352
$ws.Range("N6").Value = 0.6706928342145362
$ws.Range("O6").Value = 0.6823518567463578
$ws.Range("P6").Value = 0.6912963083250303
$ws.Range("Q6").Value = 0.8217139100174641
$ws.Range("R6").Value = 0.8217139100174641
$ws.Range("S6").Value = 0.9651392705172921
$ws.Range("T6").Value = 0.9731753977171517
$ws.Range("U6").Value = 0.9731753977171517
$ws.Range("V6").Value = 0.9731753977171517
$ws.Range("W6").Value = 0.9731753977171517
$ws.Range("X6").Value = 0.9731753977171517
$ws.Range("Y6").Value = 0.9754074930617253
$ws.Range("Z6").Value = 0.9754074930617253
$ws.Range("AA6").Value = 0.9885736230491529
$ws.Range("AB6").Value = 0.9885736230491529
$ws.Range("AC6").Value = 0.9995717111729875
$ws.Range("AD6").Value = 0.9995717111729875
$ws.Range("AE6").Value = 1
$ws.Range("AF6").Value = 1
$ws.Range("AG6").Value = 1
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 1
$ws.Range("AJ6").Value = 1
$ws.Range("AK6").Value = 1

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = 0.6302548612520517
$ws.Range("G2").Value = 5
$ws.Range("F3").Value = 0.5643644035494009
$ws.Range("D4").Value = 10
$ws.Range("F4").Value = 0.5747840530482897
$ws.Range("G4").Value = 8
$ws.Range("F5").Value = 0.5344809922170457
$ws.Range("D6").Value = 10
$ws.Range("F6").Value = 0.541384888766529
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 13
$ws.Range("F2").Value = 0.7176299835119724
$ws.Range("G2").Value = 12
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 0.7151556571962822
$ws.Range("G3").Value = 9
$ws.Range("F4").Value = 0.7055045242959332
$ws.Range("D5").Value = 15
$ws.Range("F5").Value = 0.7706968936328267
$ws.Range("G5").Value = 13
$ws.Range("F6").Value = 0.8217139100174641

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 15
$ws.Range("F2").Value = 0.8072491810088468
$ws.Range("G2").Value = 14
$ws.Range("D3").Value = 16
$ws.Range("F3").Value = 0.8244844336096276
$ws.Range("G3").Value = 13
$ws.Range("F4").Value = 0.8201613266801981
$ws.Range("D5").Value = 16
$ws.Range("F5").Value = 0.8270523783622279
$ws.Range("G5").Value = 14
$ws.Range("F6").Value = 0.8217139100174641

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 0.9152532003953405
$ws.Range("G2").Value = 17
$ws.Range("D3").Value = 19
$ws.Range("F3").Value = 0.9117097464081362
$ws.Range("G3").Value = 16
$ws.Range("F4").Value = 0.9412862438581445
$ws.Range("D5").Value = 19
$ws.Range("F5").Value = 0.908376800553647
$ws.Range("G5").Value = 17
$ws.Range("F6").Value = 0.9651392705172921

Write-Output "Applied $([string]280) cell updates across 6 sheets"
